$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")

$ws.Range("S2").Value = 53.89737328
$ws.Range("T2").Value = 62.46705563
$ws.Range("U2").Value = 72.39931747999999
$ws.Range("V2").Value = 83.91080896
$ws.Range("W2").Value = 97.25262758
$ws.Range("X2").Value = 112.7157954
$ws.Range("Y2").Value = 130.6376068
$ws.Range("Z2").Value = 140.0435145
$ws.Range("AA2").Value = 150.1266476
$ws.Range("AB2").Value = 160.9357662
$ws.Range("AC2").Value = 172.5231414
$ws.Range("AD2").Value = 184.9448075
$ws.Range("AE2").Value = 199.3705025
$ws.Range("AF2").Value = 214.9214017
$ws.Range("AG2").Value = 231.6852711
$ws.Range("AH2").Value = 249.7567222
$ws.Range("AI2").Value = 269.2377465
$ws.Range("AJ2").Value = 292.3921927
$ws.Range("AK2").Value = 317.5379213
$ws.Range("AL2").Value = 344.8461825
$ws.Range("AM2").Value = 374.5029542
$ws.Range("AN2").Value = 406.7102083
$ws.Range("AO2").Value = 443.314127
$ws.Range("AP2").Value = 483.2123985
$ws.Range("AQ2").Value = 526.7015143
$ws.Range("AR2").Value = 574.1046506
$ws.Range("AS2").Value = 625.7740692
